$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.617.66'
$ws.Range('E2').Value = '  +2.52%  '
$ws.Range('D3').Value = '1.788.66'
$ws.Range('E3').Value = '  +0.65%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '223.37'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '32.83'
$ws.Range('E8').Value = '  +7.61%  '
$ws.Range('E9').Value = '  +0.62%  '
$ws.Range('D10').Value = '0.0678'
$ws.Range('E10').Value = '  +1.65%  '
$ws.Range('E11').Value = '  +1.40%  '
$ws.Range('D12').Value = '2.045.09'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').Value = '11.16'
$ws.Range('E13').Value = '  +11.06%  '
$ws.Range('D14').Value = '1.789.50'
$ws.Range('E14').Value = '  +0.66%  '
$ws.Range('D15').Value = '34.617.11'
$ws.Range('E15').Value = '  +2.60%  '
$ws.Range('D16').Value = '0.632'
$ws.Range('E16').Value = '  +0.88%  '
$ws.Range('E17').Value = '  +2.97%  '
$ws.Range('D18').Value = '68.54'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = '253.25'
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('E20').Value = '  +4.42%  '
$ws.Range('D21').Value = '0.999'
$ws.Range('E21').Value = '  +0.05%  '
$ws.Range('D22').Value = '10.45'
$ws.Range('E22').Value = '  +1.62%  '
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('D25').Value = '159.05'
$ws.Range('E25').Value = '  -0.25%  '
$ws.Range('D26').Value = '16.32'
$ws.Range('E26').Value = '  -1.05%  '
$ws.Range('E27').Value = '  +2.02%  '
$ws.Range('E28').Value = '  -0.12%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').Value = '  -1.71%  '
$ws.Range('D31').Value = '0.0515'
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('E32').Value = '  -0.44%  '
$ws.Range('D33').Value = '3.58'
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('E34').Value = '  -0.63%  '
$ws.Range('D35').Value = '1.440.69'
$ws.Range('E35').Value = '  -3.04%  '
$ws.Range('E36').Value = '  -1.14%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '0.0189'
$ws.Range('E37').Value = '  +2.42%  '
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').Value = '0.631'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').Value = '83.09'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('E40').Value = '  +4.03%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').Value = '0.906'
$ws.Range('E42').Value = '  +2.20%  '
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('D44').Value = '0.0506'
$ws.Range('E44').Value = '  -1.08%  '
$ws.Range('D45').Value = '5.93'
$ws.Range('E45').Value = '  +3.47%  '
$ws.Range('D46').Value = '1.05'
$ws.Range('E46').Value = '  -2.23%  '
$ws.Range('D47').Value = '1.944.77'
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('D48').Value = '104.81'
$ws.Range('E48').Value = '  +7.39%  '
$ws.Range('D49').Value = '12.02'
$ws.Range('E49').Value = '  +2.17%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').Value = '49.42'
$ws.Range('E51').Value = '  -2.64%  '
